$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 573
$ws.Range("F3").Value = 251
$ws.Range("F4").Value = 27
$ws.Range("F5").Value = 733
$ws.Range("F6").Value = 359
$ws.Range("F8").Value = 144
$ws.Range("F9").Value = 240
$ws.Range("F10").Value = 211
$ws.Range("F11").Value = 5881
$ws.Range("F12").Value = 50
$ws.Range("F13").Value = 38
$ws.Range("F16").Value = 545
$ws.Range("F17").Value = 350
$ws.Range("F18").Value = 420
$ws.Range("F21").Value = 703
$ws.Range("F22").Value = 127
$ws.Range("F23").Value = 92
$ws.Range("F24").Value = 304
$ws.Range("F25").Value = 1011
$ws.Range("F26").Value = 63
$ws.Range("F27").Value = 1790
$ws.Range("F28").Value = 460

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 659
$ws.Range("F3").Value = 265
$ws.Range("F4").Value = 50
$ws.Range("F5").Value = 266
$ws.Range("F7").Value = 72
$ws.Range("F8").Value = 46

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 210

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 210
$ws.Range("F3").Value = 573
$ws.Range("F4").Value = 251
$ws.Range("F5").Value = 27
$ws.Range("F6").Value = 733
$ws.Range("F7").Value = 659
$ws.Range("F8").Value = 359
$ws.Range("F10").Value = 144
$ws.Range("F11").Value = 240
$ws.Range("F12").Value = 211
$ws.Range("F13").Value = 5881
$ws.Range("F14").Value = 50
$ws.Range("F15").Value = 38
$ws.Range("F16").Value = 265
$ws.Range("F19").Value = 545
$ws.Range("F20").Value = 350
$ws.Range("F21").Value = 420
$ws.Range("F22").Value = 50
$ws.Range("F25").Value = 266
$ws.Range("F27").Value = 72
$ws.Range("F28").Value = 703
$ws.Range("F29").Value = 46
$ws.Range("F32").Value = 127
$ws.Range("F33").Value = 92
$ws.Range("F34").Value = 304
$ws.Range("F35").Value = 1011
$ws.Range("F36").Value = 63
$ws.Range("F37").Value = 1790
$ws.Range("F38").Value = 460
